$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (MyOrders): fill in Runmode "Done" and flip Y -> N
$ws.Range("B15").Value = "Done"
$ws.Range("C15").Value = "N"

# Row 19 (Settings): fill in Runmode "Done"
$ws.Range("B19").Value = "Done"

# New row 20: AdminSearch test case, still pending (Y)
$ws.Range("A20").Value = "AdminSearch"
$ws.Range("C20").Value = "Y"

# Update selection to reflect last edited cell
$ws.Range("C21").Select()
